$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 5 data rows (rows 2-6), one per year (2019-2023) for
# MapCode "NO2". We need to insert one new "DE_LU" row above each existing
# "NO2" row (for the matching year), pushing the table from A1:G6 to A1:G11.
#
# Insert from the bottom up so the row indices of rows not yet processed stay
# stable while we work.

# New row for year 2023 (DE_LU), goes above existing row 6 (2023/NO2)
$ws.Rows.Item(6).Insert(-4142)
$ws.Rows.Item(6).ClearFormats()
$ws.Cells.Item(6,1).Value = 2023
$ws.Cells.Item(6,2).Value = "DE_LU"
$ws.Cells.Item(6,3).Value = 95.18000000000001
$ws.Cells.Item(6,4).Value = 98.02
$ws.Cells.Item(6,5).Value = 47.58
$ws.Cells.Item(6,6).Value = -500
$ws.Cells.Item(6,7).Value = 524.27

# New row for year 2022 (DE_LU), goes above existing row 5 (2022/NO2)
$ws.Rows.Item(5).Insert(-4142)
$ws.Rows.Item(5).ClearFormats()
$ws.Cells.Item(5,1).Value = 2022
$ws.Cells.Item(5,2).Value = "DE_LU"
$ws.Cells.Item(5,3).Value = 235.44
$ws.Cells.Item(5,4).Value = 208.34
$ws.Cells.Item(5,5).Value = 142.82
$ws.Cells.Item(5,6).Value = -19.04
$ws.Cells.Item(5,7).Value = 871

# New row for year 2021 (DE_LU), goes above existing row 4 (2021/NO2)
$ws.Rows.Item(4).Insert(-4142)
$ws.Rows.Item(4).ClearFormats()
$ws.Cells.Item(4,1).Value = 2021
$ws.Cells.Item(4,2).Value = "DE_LU"
$ws.Cells.Item(4,3).Value = 96.84999999999999
$ws.Cells.Item(4,4).Value = 75.48
$ws.Cells.Item(4,5).Value = 73.68000000000001
$ws.Cells.Item(4,6).Value = -69
$ws.Cells.Item(4,7).Value = 620

# New row for year 2020 (DE_LU), goes above existing row 3 (2020/NO2)
$ws.Rows.Item(3).Insert(-4142)
$ws.Rows.Item(3).ClearFormats()
$ws.Cells.Item(3,1).Value = 2020
$ws.Cells.Item(3,2).Value = "DE_LU"
$ws.Cells.Item(3,3).Value = 30.47
$ws.Cells.Item(3,4).Value = 30.99
$ws.Cells.Item(3,5).Value = 17.5
$ws.Cells.Item(3,6).Value = -83.94
$ws.Cells.Item(3,7).Value = 200.04

# New row for year 2019 (DE_LU), goes above existing row 2 (2019/NO2)
$ws.Rows.Item(2).Insert(-4142)
$ws.Rows.Item(2).ClearFormats()
$ws.Cells.Item(2,1).Value = 2019
$ws.Cells.Item(2,2).Value = "DE_LU"
$ws.Cells.Item(2,3).Value = 37.67
$ws.Cells.Item(2,4).Value = 38.06
$ws.Cells.Item(2,5).Value = 15.52
$ws.Cells.Item(2,6).Value = -90.01000000000001
$ws.Cells.Item(2,7).Value = 121.46
